$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Text = "2025-09-13 Saturday"

# Update each table cell value (row-major order, 5 columns)
$t = $d.Tables.Item(1)
$values = @(
    "13+82=",
    "36-18=",
    "53+45=",
    "90-18=",
    "66-45=",
    "64+28=",
    "70-51=",
    "95+0=",
    "34-11=",
    "13+77=",
    "48+50=",
    "94-53=",
    "92-55=",
    "52+17=",
    "27+34=",
    "31+23=",
    "34+28=",
    "50+0=",
    "80-50=",
    "97-12=",
    "82-10=",
    "53+11=",
    "77+9=",
    "20+57=",
    "41-1=",
    "39-28=",
    "61-20=",
    "67-56=",
    "63-38=",
    "27-2=",
    "27-18=",
    "21+36=",
    "2+33=",
    "64-44=",
    "21+72=",
    "16+22=",
    "44-28=",
    "59-34=",
    "85+1=",
    "89-86=",
    "64+25=",
    "53-5=",
    "48+39=",
    "16-15=",
    "26+33=",
    "22-4=",
    "50+44=",
    "65-48=",
    "43+12=",
    "86+13=",
    "8+48=",
    "50+19=",
    "45-36=",
    "74-35=",
    "69-8=",
    "23+55=",
    "23+16=",
    "17+57=",
    "83-23=",
    "21-17=",
    "45+48=",
    "51+47=",
    "33+37=",
    "53-24=",
    "82+7=",
    "87-64=",
    "54-25=",
    "88-18=",
    "18-1=",
    "68-37=",
    "37+16=",
    "26+34=",
    "92-78=",
    "74-40=",
    "98-79=",
    "31-19=",
    "29+62=",
    "62-49=",
    "46+38=",
    "88-18=",
    "4+23=",
    "10+37=",
    "73-9=",
    "4+67=",
    "54-54=",
    "14+79=",
    "77-44=",
    "10+41=",
    "28+11=",
    "82-37=",
    "51-26=",
    "45+11=",
    "70+19=",
    "10+18=",
    "43-15=",
    "24-0=",
    "25+57=",
    "94-63=",
    "46-44=",
    "5+85="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "done"
